$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.373
$ws.Range("D3").Value = -7.428999999999999
$ws.Range("E8").Value = 16.842
$ws.Range("E11").Value = 17.177
$ws.Range("B12").Value = 4.871
$ws.Range("C14").Value = -12.987
$ws.Range("E14").Value = 16.858
$ws.Range("E15").Value = 15.916
$ws.Range("E17").Value = 16.771
$ws.Range("D20").Value = -7.755000000000001
$ws.Range("D25").Value = -8.163
$ws.Range("C26").Value = -11.508
$ws.Range("E26").Value = 16.777
$ws.Range("B27").Value = 5.422
$ws.Range("D30").Value = -7.176
$ws.Range("C31").Value = -12.516
$ws.Range("B32").Value = 5.733
$ws.Range("C35").Value = -12.35
$ws.Range("B36").Value = 8.620999999999999
$ws.Range("E36").Value = 16.403
$ws.Range("C37").Value = -13.665
$ws.Range("B38").Value = 5.148000000000001
$ws.Range("D44").Value = -7.749
$ws.Range("C45").Value = -12.776
$ws.Range("B46").Value = 5.747999999999999
$ws.Range("D47").Value = -7.579000000000001
$ws.Range("C52").Value = -11.07
$ws.Range("B54").Value = 5.119000000000002
$ws.Range("B55").Value = 4.595
$ws.Range("B56").Value = 4.614
$ws.Range("C57").Value = -13.697
$ws.Range("D58").Value = -8.021000000000001
$ws.Range("E64").Value = 17.444
$ws.Range("B67").Value = 5.212000000000001
$ws.Range("B69").Value = 5.211999999999999
$ws.Range("B72").Value = 5.624
$ws.Range("D78").Value = -7.906999999999999
$ws.Range("E79").Value = 17.571
$ws.Range("C81").Value = -13.352
$ws.Range("B83").Value = 5.007000000000001
$ws.Range("C83").Value = -13.772
$ws.Range("D84").Value = -8.125999999999999
$ws.Range("B86").Value = 5.037999999999999
$ws.Range("D89").Value = -7.363
$ws.Range("E89").Value = 17.235
$ws.Range("B91").Value = 5.228
$ws.Range("D91").Value = -6.857000000000001
$ws.Range("D92").Value = -6.884
$ws.Range("B93").Value = 5.77
$ws.Range("D96").Value = -7.540999999999999
$ws.Range("B99").Value = 5.811
$ws.Range("C100").Value = -12.662
$ws.Range("C102").Value = -13.228
$ws.Range("D102").Value = -7.529999999999999
